# Updates cryptos list prices / 1h-volume percentages (and swaps the
# Aave / OKB rows) to match the latest scrape, per commit
# "Updated cryptos list on Sat Sep 28 04:50:25 UTC 2024 with GitHub Actions".
#
# Cells in column D that look like plain numbers are written with a
# leading apostrophe so Excel keeps them as text (matching the source
# data, which stores prices such as "1.00" / "42.58" as strings rather
# than numeric values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.045.20"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "2.697.14"
$ws.Range("E3").Value = "  +1.93%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'612.16"
$ws.Range("E5").Value = "  +1.41%  "
$ws.Range("D6").Value = "'158.55"
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  +4.21%  "
$ws.Range("D10").Value = "'6.04"
$ws.Range("E10").Value = "  +4.08%  "
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  +9.54%  "
$ws.Range("D14").Value = "'30.15"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").Value = "3.183.65"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").Value = "65.915.19"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "2.695.93"
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "'12.78"
$ws.Range("E18").Value = "  +0.50%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("E20").Value = "  +5.64%  "
$ws.Range("D21").Value = "'359.32"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").Value = "'71.42"
$ws.Range("E22").Value = "  +3.11%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("E24").Value = "  +17.13%  "
$ws.Range("D25").Value = "'10.01"
$ws.Range("E25").Value = "  +5.87%  "
$ws.Range("E26").Value = "  -3.61%  "
$ws.Range("D27").Value = "'1.67"
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("D29").Value = "'8.29"
$ws.Range("E29").Value = "  -0.21%  "
$ws.Range("E30").Value = "  +1.15%  "
$ws.Range("D32").Value = "'534.27"
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("D33").Value = "'1.79"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  +4.55%  "
$ws.Range("D35").Value = "'5.49"
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("D36").Value = "'0.434"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "'20.78"
$ws.Range("E37").Value = "  +0.72%  "
$ws.Range("D38").Value = "'162.34"
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("D39").Value = "'1.99"
$ws.Range("E39").Value = "  -1.32%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'42.58"
$ws.Range("E42").Value = "  -0.37%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'168.20"
$ws.Range("E43").Value = "  +0.84%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "'0.0636"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").Value = "'23.83"
$ws.Range("E46").Value = "  +2.34%  "
$ws.Range("D47").Value = "'2.31"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").Value = "'20.95"
$ws.Range("E50").Value = "  +6.07%  "
$ws.Range("E51").Value = "  +1.28%  "
